$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)
$ws.Rows("16:19").Insert()
$ws.Rows("14:14").Copy()
$ws.Rows("16:19").PasteSpecial(-4122)
